$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-upload of the attendance log: every logged date (column C, rows 2-75)
# moves from 2025-07-09 (serial 45847) to 2025-09-07 (serial 45907).
$ws.Range("C2:C75").Value = 45907

# Copy the (now refreshed) date formatting down onto the header cell C1 as
# well, matching the column's date style.
$ws.Range("C2").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Clear the explicit 15.75pt row heights left over from the previous Excel
# build so rows fall back to the sheet's (new) default height.
$ws.Range("A1:F75").EntireRow.AutoFit()

# Scroll back to the top and select the header row, like the re-saved file.
$ws.Rows(1).Select()
